$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row -- append "(English)" hints to the existing headers and add
#    a brand-new "Phone number" column (J).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,2).Value  = "Фамилия (Surname)"
$ws.Cells.Item(1,3).Value  = "Имя (Name)"
$ws.Cells.Item(1,4).Value  = "Курс (Year)"
$ws.Cells.Item(1,5).Value  = "Адрес электронной почты (Email)"
$ws.Cells.Item(1,6).Value  = "Адрес проживания (Address)"
$ws.Cells.Item(1,7).Value  = "Пожелания (Wishes)"
$ws.Cells.Item(1,8).Value  = "Аккаунт в социальной сети (vk ID)"
$ws.Cells.Item(1,9).Value  = "Индекс (Index)"

# New header cell for the phone-number column, copy formatting from the
# neighbouring header first so it picks up the same look & feel.
$ws.Cells.Item(1,9).Copy()
$ws.Cells.Item(1,10).PasteSpecial(-4122)
$ws.Cells.Item(1,10).Value = "Номер телефона (Phone number)"

# ---------------------------------------------------------------------------
# 2. Row 2 -- Valeriya Terova's submission: name correction + new phone cell.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = 44170.6170574537
$ws.Cells.Item(2,3).Value = "Валерия"
$ws.Cells.Item(2,7).Value = "Спать"
$ws.Cells.Item(2,9).Value = 45678.0

$ws.Cells.Item(2,9).Copy()
$ws.Cells.Item(2,10).PasteSpecial(-4122)
$ws.Cells.Item(2,10).Value = "+7-921-873-6059"

# ---------------------------------------------------------------------------
# 3. Row 3 -- replaced by a new respondent (Platon Platonov).
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,1).Value = 44170.61825546296
$ws.Cells.Item(3,2).Value = "Платонов"
$ws.Cells.Item(3,3).Value = "Платон"
$ws.Cells.Item(3,4).Value = "2 магистратуры (master)"
$ws.Cells.Item(3,7).Value = "Есть"
$ws.Cells.Item(3,9).Value = 456789.0

# ---------------------------------------------------------------------------
# 4. Row 4 -- brand-new respondent (Stepan Prepod).
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(4,1).PasteSpecial(-4122)
$ws.Cells.Item(4,1).Value = 44170.61911962963

$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(4,2).PasteSpecial(-4122)
$ws.Cells.Item(4,2).Value = "Препод"

$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(4,3).PasteSpecial(-4122)
$ws.Cells.Item(4,3).Value = "Степан"

$ws.Cells.Item(3,5).Copy()
$ws.Cells.Item(4,5).PasteSpecial(-4122)
$ws.Cells.Item(4,5).Value = "terovaleriya@ya.ru"

$ws.Cells.Item(3,6).Copy()
$ws.Cells.Item(4,6).PasteSpecial(-4122)
$ws.Cells.Item(4,6).Value = "14 линия"

$ws.Cells.Item(3,7).Copy()
$ws.Cells.Item(4,7).PasteSpecial(-4122)
$ws.Cells.Item(4,7).Value = "Умных студентов"

$ws.Cells.Item(3,9).Copy()
$ws.Cells.Item(4,9).PasteSpecial(-4122)
$ws.Cells.Item(4,9).Value = 3456765.0

# ---------------------------------------------------------------------------
# 5. Widen the formatted column range so it covers the new column J too.
# ---------------------------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 20.67

$excel.CutCopyMode = $false
